# feat: update GenomeEntity Resource
#
# The "snp" (Single Nucleotide Polymorphism) block gains two new rows
# (LOINC#48000-4 / CodeableConcept and LOINC#LP232001-0 / Quantity),
# which also re-orders the row that held LOINC#74019-1 to the top of
# that block. Inserting the two rows pushes the whole downstream
# "variant-annotation" + "vital-status" block down by two rows intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: reorder the existing 3 "snp" rows (85-87) -------------------
# Before: 85=LOINC#93348-1(boolean) 86=UMLS#C0449889(boolean) 87=LOINC#74019-1(CodeableConcept)
# After:  85=LOINC#74019-1(CodeableConcept) 86=LOINC#93348-1(boolean) 87=UMLS#C0449889(boolean)
$ws.Range("E85").Value = "LOINC#74019-1"
$ws.Range("H85").Value = "CodeableConcept"

$ws.Range("E86").Value = "LOINC#93348-1"
$ws.Range("H86").Value = "boolean"

$ws.Range("E87").Value = "Unified Medical Language System#C0449889"
$ws.Range("H87").Value = "boolean"

# --- Step 2: insert two new rows at 88:89, carrying formatting down ------
$ws.Rows("88:89").Insert()

# Copy the formatting (style) of the row directly above so the new rows
# match the rest of the table (border/alignment/font) instead of picking
# up the default style.
$ws.Range("A87:K87").Copy()
$ws.Range("A88:K89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: populate the two new "snp" rows ------------------------------
$ws.Range("A88").Value = ""
$ws.Range("B88").Value = "Single Nucleotide Polymorphism"
$ws.Range("C88").Value = ""
$ws.Range("D88").Value = ""
$ws.Range("E88").Value = "LOINC#48000-4"
$ws.Range("F88").Value = ""
$ws.Range("G88").Value = ""
$ws.Range("H88").Value = "CodeableConcept"
$ws.Range("I88").Value = "optional"
$ws.Range("J88").Value = ""
$ws.Range("K88").Value = ""

$ws.Range("A89").Value = ""
$ws.Range("B89").Value = "Single Nucleotide Polymorphism"
$ws.Range("C89").Value = ""
$ws.Range("D89").Value = ""
$ws.Range("E89").Value = "LOINC#LP232001-0"
$ws.Range("F89").Value = ""
$ws.Range("G89").Value = ""
$ws.Range("H89").Value = "Quantity"
$ws.Range("I89").Value = "optional"
$ws.Range("J89").Value = ""
$ws.Range("K89").Value = ""
